$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("L2").Value = 5801
$ws.Range("L3").Value = 6318
$ws.Range("C4").Value = 1876
$ws.Range("L4").Value = 1565
$ws.Range("L5").Value = 378
$ws.Range("L6").Value = 5188
$ws.Range("C7").Value = 28420
$ws.Range("L7").Value = 19250

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("L2").Value = 161
$ws.Range("L7").Value = 616
$ws.Range("L8").Value = 1267
$ws.Range("L11").Value = 319
$ws.Range("L16").Value = 41
$ws.Range("L19").Value = 523
$ws.Range("L23").Value = 211
$ws.Range("L29").Value = 1088
$ws.Range("L33").Value = 874
$ws.Range("L37").Value = 735
$ws.Range("L39").Value = 13
$ws.Range("L42").Value = 621
$ws.Range("L50").Value = 95
$ws.Range("L51").Value = 241
$ws.Range("L52").Value = 401
$ws.Range("L53").Value = 213
$ws.Range("L54").Value = 421
$ws.Range("L55").Value = 200
$ws.Range("C63").Value = 300
$ws.Range("K63").Value = 177
$ws.Range("L63").Value = 56
$ws.Range("K64").Value = 167
$ws.Range("L65").Value = 373
$ws.Range("L67").Value = 667
$ws.Range("L74").Value = 16
$ws.Range("L76").Value = 293
$ws.Range("L78").Value = 246
$ws.Range("L83").Value = 423
$ws.Range("L85").Value = 954
$ws.Range("L88").Value = 202
$ws.Range("L91").Value = 258
$ws.Range("L96").Value = 220
$ws.Range("L98").Value = 103
$ws.Range("L99").Value = 336
$ws.Range("C101").Value = 28420
$ws.Range("L101").Value = 19250

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("L6").Value = 67
$ws.Range("L7").Value = 220

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("L2").Value = 211
$ws.Range("L4").Value = 45
$ws.Range("L7").Value = 616

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("L3").Value = 96
$ws.Range("L6").Value = 78
$ws.Range("L7").Value = 319

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("L2").Value = 289
$ws.Range("L5").Value = 21
$ws.Range("L6").Value = 198
$ws.Range("L7").Value = 954

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("L2").Value = 127
$ws.Range("L7").Value = 401

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("L2").Value = 63
$ws.Range("L7").Value = 213

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("L2").Value = 377
$ws.Range("L3").Value = 448
$ws.Range("L7").Value = 1267

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("L2").Value = 134
$ws.Range("L6").Value = 93
$ws.Range("L7").Value = 423

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("L5").Value = 22
$ws.Range("L6").Value = 250
$ws.Range("L7").Value = 874

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("L2").Value = 221
$ws.Range("L3").Value = 257
$ws.Range("L7").Value = 735

$ws = $wb.Worksheets.Item('New City')
$ws.Range("L2").Value = 135
$ws.Range("L3").Value = 123
$ws.Range("L7").Value = 373

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("L3").Value = 137
$ws.Range("L7").Value = 336

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("L3").Value = 258
$ws.Range("L7").Value = 667

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("L6").Value = 206
$ws.Range("L7").Value = 421

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("L2").Value = 322
$ws.Range("L3").Value = 421
$ws.Range("L7").Value = 1088

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("L3").Value = 163
$ws.Range("L7").Value = 523

$ws = $wb.Worksheets.Item('River North')
$ws.Range("L2").Value = 61
$ws.Range("L7").Value = 293

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("L2").Value = 170
$ws.Range("L7").Value = 621

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("L4").Value = 27
$ws.Range("L6").Value = 70
$ws.Range("L7").Value = 246

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("L4").Value = 18
$ws.Range("L7").Value = 200

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("L2").Value = 53
$ws.Range("L7").Value = 211

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("L6").Value = 34
$ws.Range("L7").Value = 258

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("K5").Value = 5
$ws.Range("K7").Value = 167

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("L4").Value = 11
$ws.Range("L7").Value = 103

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("L6").Value = 24
$ws.Range("L7").Value = 95

$ws = $wb.Worksheets.Item('Greektown')
$ws.Range("L2").Value = 6
$ws.Range("L6").Value = 13

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("L2").Value = 53
$ws.Range("L7").Value = 161

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("L3").Value = 69
$ws.Range("L7").Value = 202

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("L3").Value = 25
$ws.Range("L4").Value = 69

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("L2").Value = 72
$ws.Range("L7").Value = 241

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("L2").Value = 7
$ws.Range("L7").Value = 41

$ws = $wb.Worksheets.Item('Printers Row')
$ws.Range("L3").Value = 6
$ws.Range("L7").Value = 16
